# Updated cryptos list on Wed Feb 22 14:25:25 UTC 2023 with GitHub Actions
#
# Applies the latest price/volume(1h) snapshot to the crypto tracker sheet.
# Columns: A=rank(unchanged) B=Coin C=Link D=Price E=Volume(1h)
# Rows 10/11 and 36/37 swapped rank order (OKB<->Polygon, ImmutableX<->VeChain),
# so Coin/Link/Price/Volume are rewritten for those rows too.
#
# D-column price cells are free-form text (e.g. "24.090.83" with thousand
# dots, or plain decimals like "308.74"). Excel auto-converts plain
# decimal-looking text to a Number on assignment, so for any new D value
# that parses as a number we force the cell to Text ("@") first to keep it
# a literal string, matching the original inlineStr cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.088.02"
$ws.Range("E2").Value = "  -2.63%  "
$ws.Range("D3").Value = "1.637.32"
$ws.Range("E3").Value = "  -2.61%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.89"
$ws.Range("E5").Value = "  -2.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3954"
$ws.Range("E7").Value = "  +0.74%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3853"
$ws.Range("E8").Value = "  -2.75%  "
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "49.78"
$ws.Range("E10").Value = "  -3.80%  "
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.362"
$ws.Range("E11").Value = "  -2.82%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08548"
$ws.Range("E12").Value = "  -1.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.57"
$ws.Range("E13").Value = "  -6.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.071"
$ws.Range("E14").Value = "  -3.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001282"
$ws.Range("E15").Value = "  -2.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.485"
$ws.Range("E16").Value = "  -3.55%  "
$ws.Range("D17").Value = "1.644.60"
$ws.Range("E17").Value = "  -2.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.81"
$ws.Range("E18").Value = "  +0.63%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06921"
$ws.Range("E19").Value = "  -2.45%  "
$ws.Range("E20").Value = "  +0.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.910"
$ws.Range("E21").Value = "  -2.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.003"
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.61"
$ws.Range("E23").Value = "  -2.22%  "
$ws.Range("D24").Value = "24.090.99"
$ws.Range("E24").Value = "  -2.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.407"
$ws.Range("E25").Value = "  +2.46%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.850"
$ws.Range("E26").Value = "  +0.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.19"
$ws.Range("E27").Value = "  -5.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "157.65"
$ws.Range("E28").Value = "  -3.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "139.76"
$ws.Range("E29").Value = "  -4.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.098"
$ws.Range("E30").Value = "  +2.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.247"
$ws.Range("E31").Value = "  -10.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.492"
$ws.Range("E32").Value = "  +4.87%  "
$ws.Range("D33").Value = "1.831.28"
$ws.Range("E33").Value = "  -2.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08083"
$ws.Range("E34").Value = "  -3.70%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.745"
$ws.Range("E35").Value = "  -3.05%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9697"
$ws.Range("E36").Value = "  -3.08%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02898"
$ws.Range("E37").Value = "  -5.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2687"
$ws.Range("E38").Value = "  -3.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.09241"
$ws.Range("E39").Value = "  -2.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.30"
$ws.Range("E40").Value = "  -2.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.431"
$ws.Range("E41").Value = "  -7.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7507"
$ws.Range("E42").Value = "  -5.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.04"
$ws.Range("E43").Value = "  -3.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.07"
$ws.Range("E44").Value = "  -3.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6896"
$ws.Range("E45").Value = "  -3.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.457"
$ws.Range("E46").Value = "  -4.19%  "
$ws.Range("E47").Value = "  -2.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.003"
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08341"
$ws.Range("E49").Value = "  -3.67%  "
$ws.Range("E50").Value = "  -5.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "133.47"
$ws.Range("E51").Value = "  -3.13%  "
